$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells receiving numeric-looking text values need NumberFormat forced to Text
# first, otherwise Excel auto-converts the string into a floating point number.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.842.04"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "1.640.41"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "216.02"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "0.5063"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "0.06437"
$ws.Range("E9").Value = "  +1.35%  "
$ws.Range("D10").Value = "20.49"
$ws.Range("E10").Value = "  +5.22%  "
$ws.Range("D11").Value = "0.07790"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "4.276"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "1.646.57"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").Value = "1.866.61"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").Value = "0.5633"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "0.0₅7656"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "63.37"
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "25.863.00"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "4.388"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "193.14"
$ws.Range("E21").Value = "  -0.97%  "
$ws.Range("D22").Value = "9.935"
$ws.Range("E22").Value = "  +0.71%  "
$ws.Range("D23").Value = "6.156"
$ws.Range("E23").Value = "  +1.69%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "1.805"
$ws.Range("E25").Value = "  -4.77%  "
$ws.Range("D26").Value = "141.33"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "0.1237"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").Value = "6.817"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Value = "1.246"
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("D31").Value = "0.04958"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").Value = "3.295"
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "3.240"
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").Value = "1.574"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").Value = "2.387"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").Value = "0.9056"
$ws.Range("E36").Value = "  +0.97%  "
$ws.Range("D37").Value = "0.5570"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("D38").Value = "1.133.91"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").Value = "2.548"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "0.01569"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.8046"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.478"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("D44").Value = "98.94"
$ws.Range("E44").Value = "  +1.44%  "
$ws.Range("D45").Value = "1.777.61"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "0.0₈112"
$ws.Range("E46").Value = "  -4.56%  "
$ws.Range("D47").Value = "55.76"
$ws.Range("E47").Value = "  +1.81%  "
$ws.Range("D48").Value = "0.4276"
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").Value = "7.759"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").Value = "0.05034"
$ws.Range("D51").Value = "0.9996"
$ws.Range("E51").Value = "  -0.52%  "

# Restore default (Normal) style on the cells we touched so no stray
# number-format style lingers on them.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
